# Applies the brochure_ru.pptx slide-2 layout tweaks:
#   - TextBox 23 (id 34, {{PASSWORD}}):  reposition/resize + set run font to "Circe"
#   - TextBox 34 (id 35, {{QR_WIFI}}):   reposition/resize + turn word-wrap on (wrap="square")
#   - TextBox 20 (id 36, "Код безопасен..."): reposition only
#
# NOTE on the numeric literals below: Shape.Left/.Top/.Width/.Height are exposed
# as single-precision (float32) points in this object model, then multiplied by
# 12700 to land back on EMU when the OOXML is serialized. A "clean" EMU/12700
# division loses precision in that round-trip and can land 1 EMU short, so the
# literals here are the nearest float32 points value that reliably reproduces
# the exact target EMU offsets/extents from the diff.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

function Get-ShapeById($slide, $id) {
    foreach ($sh in $slide.Shapes) {
        if ($sh.Id -eq $id) {
            return $sh
        }
    }
    return $null
}

# --- Shape id=34, "TextBox 23" -> {{PASSWORD}} ---
$shPassword = Get-ShapeById $s 34
$shPassword.Left = 69.99811553955078
$shPassword.Top = 280.71820068359375
$shPassword.Width = 185.0732421875
$shPassword.Height = 31.201732635498047
$shPassword.TextFrame.TextRange.Font.Name = "Circe"

# --- Shape id=35, "TextBox 34" -> {{QR_WIFI}} ---
$shQr = Get-ShapeById $s 35
$shQr.Left = 90.01559448242188
$shQr.Top = 386.1419982910156
$shQr.Width = 102.98441314697266
$shQr.Height = 29.081260681152344
$shQr.TextFrame.WordWrap = -1

# --- Shape id=36, "TextBox 20" -> "Код безопасен..." ---
$shCode = Get-ShapeById $s 36
$shCode.Left = 48.721893310546875
$shCode.Top = 326.8350524902344
